$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "Aluno: mariano joão" "Aluno: chiquin"
Replace-Text "Turma: 6ºB" "Turma: 6ºA"
Replace-Text "Data: 18 de junho de 2025" "Data: 23 de junho de 2025"
Replace-Text "Art. 54 – {{descricaoInciso}}" "Art. 54 – 14"
Replace-Text "Esta medida acarreta perda de sua nota disciplinar em -0.70 pontos, enquadrando-se no comportamento Bom." "Esta medida acarreta perda de sua nota disciplinar em -0.50 pontos, enquadrando-se no comportamento Excepcional."
Replace-Text "Nota Anterior: 8.32" "Nota Anterior: 10.00"
Replace-Text "Nota Atual: 7.62" "Nota Atual: 9.50"
Replace-Text "ok, alterando observação, lhe encaminhando para devidas providencias" "teste"
Replace-Text "Cruzeiro do Sul – AC, 18 de junho de 2025" "Cruzeiro do Sul – AC, 23 de junho de 2025"
